# Update the single data row (row 2) on the active sheet with the next
# scraped firm contact ("Aoife Bradley" / Byrne Wallace / Ireland), replacing
# the previous contact ("Junghwan Bae" / Bae Kim And Lee / Korea (South)).
#
# Columns (per header row 1): A=Name, B=Role, C=Firm, D=Country,
# E=Nationality, F=Practice Area, G=Email, H=Phone

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Aoife Bradley"
$ws.Range("B2").Value = "Partner"
$ws.Range("C2").Value = "Byrne Wallace"
$ws.Range("D2").Value = "Ireland"
$ws.Range("E2").Value = "-----"
$ws.Range("F2").Value = "-----"
$ws.Range("G2").Value = "abradley@byrnewallaceshields.com"

# The phone number is a purely-numeric string that must stay a text value
# (matching the shared-string cell type used by the source data) instead of
# being auto-coerced into a numeric cell. A leading apostrophe forces Excel
# to store it as text while keeping the literal digits as the cell content.
$ws.Range("H2").Value = "'35316371583"
